$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some new Price values look numeric (e.g. "412.28") but must stay stored
# as text, matching the rest of the sheet which keeps prices as inline
# strings (some contain multiple "." thousand separators, e.g. "62.506.81",
# so the whole Price column is textual). Force a Text format on those
# specific cells before writing so Excel does not coerce them into numeric
# values / alter their precision.
$textCellAddresses = @("D5", "D6", "D10", "D11", "D12", "D15", "D16", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $textCellAddresses) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '62.506.81'
$ws.Range("E2").Value = '  +5.60%  '
# Row 3
$ws.Range("D3").Value = '3.455.80'
$ws.Range("E3").Value = '  +3.04%  '
# Row 4
$ws.Range("E4").Value = '  +0.14%  '
# Row 5
$ws.Range("D5").Value = '412.28'
$ws.Range("E5").Value = '  -0.20%  '
# Row 6
$ws.Range("D6").Value = '128.45'
$ws.Range("E6").Value = '  +14.97%  '
# Row 7
$ws.Range("D7").Value = '3.449.47'
$ws.Range("E7").Value = '  +3.03%  '
# Row 8
$ws.Range("E8").Value = '  +0.81%  '
# Row 9
$ws.Range("E9").Value = '  +0.08%  '
# Row 10
$ws.Range("D10").Value = '0.687'
$ws.Range("E10").Value = '  +7.86%  '
# Row 11
$ws.Range("D11").Value = '0.125'
$ws.Range("E11").Value = '  +26.68%  '
# Row 12
$ws.Range("D12").Value = '43.54'
$ws.Range("E12").Value = '  +8.31%  '
# Row 13
$ws.Range("E13").Value = '  -0.15%  '
# Row 14
$ws.Range("D14").Value = '4.009.08'
$ws.Range("E14").Value = '  +3.04%  '
# Row 15
$ws.Range("D15").Value = '8.72'
$ws.Range("E15").Value = '  +3.16%  '
# Row 16
$ws.Range("D16").Value = '20.15'
$ws.Range("E16").Value = '  +3.53%  '
# Row 17
$ws.Range("D17").Value = '3.420.57'
$ws.Range("E17").Value = '  +1.90%  '
# Row 18
$ws.Range("D18").Value = '62.522.46'
$ws.Range("E18").Value = '  +5.81%  '
# Row 19
$ws.Range("E19").Value = '  -0.23%  '
# Row 20
$ws.Range("D20").Value = '10.99'
$ws.Range("E20").Value = '  +0.83%  '
# Row 21
$ws.Range("D21").Value = '0.0000133'
$ws.Range("E21").Value = '  +21.05%  '
# Row 22
$ws.Range("D22").Value = '3.36'
$ws.Range("E22").Value = '  -0.19%  '
# Row 23
$ws.Range("D23").Value = '13.18'
$ws.Range("E23").Value = '  +0.89%  '
# Row 24
$ws.Range("D24").Value = '81.64'
$ws.Range("E24").Value = '  +8.37%  '
# Row 25
$ws.Range("D25").Value = '311.14'
$ws.Range("E25").Value = '  +2.47%  '
# Row 26
$ws.Range("E26").Value = '  -1.33%  '
# Row 27
$ws.Range("D27").Value = '30.20'
$ws.Range("E27").Value = '  +5.09%  '
# Row 28
$ws.Range("D28").Value = '7.80'
$ws.Range("E28").Value = '  +4.67%  '
# Row 29
$ws.Range("D29").Value = '8.06'
$ws.Range("E29").Value = '  +0.59%  '
# Row 30
$ws.Range("E30").Value = '  +6.76%  '
# Row 31
$ws.Range("B31").Value = 'LEO'
$ws.Range("C31").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D31").Value = '4.37'
$ws.Range("E31").Value = '  -2.43%  '
# Row 32
$ws.Range("B32").Value = 'Kaspa'
$ws.Range("C32").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D32").Value = '0.177'
$ws.Range("E32").Value = '  +3.48%  '
# Row 33
$ws.Range("B33").Value = 'Cosmos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D33").Value = '12.17'
$ws.Range("E33").Value = '  +5.21%  '
# Row 34
$ws.Range("B34").Value = 'InjectiveProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D34").Value = '44.34'
$ws.Range("E34").Value = '  +11.21%  '
# Row 35
$ws.Range("B35").Value = 'Toncoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D35").Value = '2.61'
$ws.Range("E35").Value = '  +22.46%  '
# Row 36
$ws.Range("E36").Value = '  +0.12%  '
# Row 37
$ws.Range("D37").Value = '0.0494'
$ws.Range("E37").Value = '  -8.26%  '
# Row 38
$ws.Range("D38").Value = '52.70'
$ws.Range("E38").Value = '  +1.14%  '
# Row 39
$ws.Range("D39").Value = '3.56'
$ws.Range("E39").Value = '  +1.15%  '
# Row 40
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").Value = '  -0.18%  '
# Row 41
$ws.Range("D41").Value = '3.02'
$ws.Range("E41").Value = '  -5.23%  '
# Row 42
$ws.Range("D42").Value = '2.00'
$ws.Range("E42").Value = '  +3.86%  '
# Row 43
$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D43").Value = '138.10'
$ws.Range("E43").Value = '  +0.22%  '
# Row 44
$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").Value = '0.125'
$ws.Range("E44").Value = '  +1.75%  '
# Row 45
$ws.Range("B45").Value = 'Celestia'
$ws.Range("C45").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D45").Value = '17.84'
$ws.Range("E45").Value = '  +4.41%  '
# Row 46
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").Value = '3.99'
$ws.Range("E46").Value = '  -0.38%  '
# Row 47
$ws.Range("B47").Value = 'TheGraph'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D47").Value = '0.288'
$ws.Range("E47").Value = '  +2.48%  '
# Row 48
$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").Value = '2.27'
$ws.Range("E48").Value = '  +0.84%  '
# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '22.44'
$ws.Range("E49").Value = '  -0.83%  '
# Row 50
$ws.Range("D50").Value = '2.226.58'
$ws.Range("E50").Value = '  +0.67%  '
# Row 51
$ws.Range("B51").Value = 'ApeXProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D51").Value = '2.40'
$ws.Range("E51").Value = '  +0.14%  '

# Restore the default style on the forced-text cells (keeps them as text,
# drops the temporary Text number format so no extra style index lingers).
foreach ($addr in $textCellAddresses) {
    $ws.Range($addr).Style = "Normal"
}
